# This script applies a row-wise permutation of the "Fecha" (D), "Volumen" (J),
# "Precio minimo" (K), "Precio maximo" (L), "Precio promedio ponderado" (M) and
# "Precio $/Kg" (P) values across data rows 2-37 of the active worksheet.
# All other columns (A, B, C, E-I, N, O, Q, R) remain untouched since they are
# identical for every row (static reference data for this market/product).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (the row whose D/J/K/L/M/P values
# should be copied into the destination row).
$rowMap = @{
    2  = 16
    3  = 19
    4  = 10
    5  = 37
    6  = 36
    7  = 21
    8  = 14
    9  = 32
    10 = 23
    11 = 25
    12 = 9
    13 = 33
    14 = 24
    15 = 29
    16 = 13
    17 = 20
    18 = 35
    19 = 3
    20 = 7
    21 = 18
    22 = 11
    23 = 12
    24 = 6
    25 = 17
    26 = 30
    27 = 4
    28 = 5
    29 = 22
    30 = 34
    31 = 8
    32 = 26
    33 = 15
    34 = 27
    35 = 31
    36 = 28
    37 = 2
}

# Columns whose values get permuted between rows.
$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the current ("before") values for the relevant columns/rows first,
# since we will overwrite cells in place and rows reference each other.
# Use .Value2 for reads (plain numeric/string value, no formatting wrapper).
$snapshot = @{}
foreach ($row in 2..37) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Apply the permutation using the captured snapshot so that every destination
# row gets the original (pre-edit) values of its mapped source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcValues[$col]
    }
}
